# B6-PowerPoint.pptx edit
#
# 1) Re-style the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style {D638B6A0-1DDC-4F20-8031-556D25A37FCB} to the
#    built-in table style {4194D932-226A-41F9-9EFA-62EA12EF30D3}.
# 2) Swap the presentation's theme palette ("Integral" / Red Violet)
#    for the Office default palette ("Office Theme" / Office colors)
#    by rewriting the 12 theme colour-scheme slots in place.

$p = $ppt.ActivePresentation

$newTableStyle = "{4194D932-226A-41F9-9EFA-62EA12EF30D3}"
$tableSlides = @(14, 15, 16)
foreach ($slideIndex in $tableSlides) {
    $slide = $p.Slides.Item($slideIndex)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# Office theme colour scheme (RGB ints use the VBA RGB() byte order,
# i.e. R + G*256 + B*65536) in the standard 12-slot theme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
